$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as text first so values like "1.000" / "0.9993"
# are stored verbatim instead of being auto-parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.447.10"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").Value = "1.815.36"
$ws.Range("E3").Value = "  +5.75%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "343.31"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.3813"
$ws.Range("E7").Value = "  +3.23%  "
$ws.Range("D8").Value = "0.3498"
$ws.Range("E8").Value = "  +4.66%  "
$ws.Range("D9").Value = "48.95"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "1.234"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +3.77%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "22.15"
$ws.Range("E13").Value = "  +10.23%  "
$ws.Range("D14").Value = "6.611"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").Value = "1.814.06"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "7.236"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").Value = "0.00001120"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "0.06715"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "86.39"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "17.61"
$ws.Range("E21").Value = "  +7.37%  "
$ws.Range("D22").Value = "6.570"
$ws.Range("E22").Value = "  +7.95%  "
$ws.Range("D23").Value = "13.24"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "27.444.89"
$ws.Range("E24").Value = "  +5.21%  "
$ws.Range("D25").Value = "2.468"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "2.665"
$ws.Range("E26").Value = "  +8.54%  "
$ws.Range("D27").Value = "22.05"
$ws.Range("E27").Value = "  +14.61%  "
$ws.Range("D28").Value = "1.471"
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("D29").Value = "153.95"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").Value = "2.016.23"
$ws.Range("E30").Value = "  +5.46%  "
$ws.Range("D31").Value = "135.99"
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("D32").Value = "6.342"
$ws.Range("E32").Value = "  +6.42%  "
$ws.Range("D33").Value = "4.039"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "13.96"
$ws.Range("E34").Value = "  +8.08%  "
$ws.Range("D35").Value = "0.08786"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "1.691"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "5.634"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").Value = "0.6959"
$ws.Range("E38").Value = "  +12.65%  "
$ws.Range("D39").Value = "0.2275"
$ws.Range("E39").Value = "  +6.21%  "
$ws.Range("D40").Value = "0.02407"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").Value = "0.06493"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("D42").Value = "8.922"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("E43").Value = "  +5.20%  "
$ws.Range("D44").Value = "14.71"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "0.6531"
$ws.Range("E45").Value = "  +10.56%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "4.018"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").Value = "2.177"
$ws.Range("E48").Value = "  +7.79%  "
$ws.Range("D49").Value = "133.20"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("D50").Value = "0.07326"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "80.65"
$ws.Range("E51").Value = "  +4.77%  "

# Remove the temporary text format so the cells fall back to the default
# (unstyled) cell format, matching the original workbook.
$ws.Range("D2:D51").ClearFormats()
